$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old "SBS Main Indicators..." row (row 51),
# pushing the existing source/citation rows down by one.
$ws.Rows(51).Insert()

# After the insert:
#   A50 = "Source:"                         (unchanged)
#   A51 = "" (new, blank)                   (unchanged - already correct)
#   A52 = "SBS Main Indicators..."          (unchanged - already correct)
#   A53 = "http://epp.eurostat..." + hyperlink, style HyperLink
#   A54 = "" (the row that used to be the blank A53)
#   A57 = "SME Performance Review EU"
#   A58 = "SME Performance Review EU, \"SBA Fact sheet\"..." (long citation)

# Move the URL text down into A54 (its target row) and blank out A53,
# removing the hyperlink and its special styling along the way.
$ws.Range("A53").Hyperlinks.Delete()
$ws.Range("A54").Value = $ws.Range("A53").Value()
$ws.Range("A54").Style = "source"
$ws.Range("A53").Value = ""
$ws.Range("A53").Style = "source"

# Shorten the long citation on the last row down to just the short title,
# matching the row above it.
$ws.Range("A58").Value = "SME Performance Review EU"
